$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts the inline string
# into a numeric cell (losing the original text formatting, e.g. "505.74").
$numericLookingCells = @("D5", "D6", "D10", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D28", "D29", "D32", "D34", "D35", "D36", "D38", "D39", "D40", "D43", "D44", "D45", "D47", "D49", "D50", "D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "56.501.96"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.377.14"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "505.74"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "130.71"
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("D9").Value = "2.389.60"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").Value = "0.0987"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  +5.65%  "
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").Value = "2.799.66"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "56.456.12"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "21.66"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "2.385.28"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "10.07"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").Value = "4.06"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "309.67"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "6.30"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "66.36"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").Value = "7.27"
$ws.Range("E28").Value = "  -2.77%  "
$ws.Range("D29").Value = "173.69"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").Value = "0.0₃0714"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").Value = "5.84"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "1.08"
$ws.Range("E35").Value = "  -4.23%  "
$ws.Range("D36").Value = "17.69"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").Value = "3.71"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("D39").Value = "0.824"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "36.52"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("E41").Value = "  -3.12%  "
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "4.96"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "128.35"
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").Value = "0.566"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").Value = "240.59"
$ws.Range("E47").Value = "  -4.55%  "
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").Value = "0.0208"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").Value = "17.07"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").Value = "0.952"
$ws.Range("E51").Value = "  -0.18%  "
